$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-30 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-31 Sunday", 2)

$d.Content.Find.Execute("226×5=", $true, $false, $false, $false, $false, $true, 1, $false, "804×8=", 2)
$d.Content.Find.Execute("659×7=", $true, $false, $false, $false, $false, $true, 1, $false, "107×6=", 2)
$d.Content.Find.Execute("187×8=", $true, $false, $false, $false, $false, $true, 1, $false, "111×4=", 2)
$d.Content.Find.Execute("643×6=", $true, $false, $false, $false, $false, $true, 1, $false, "805×3=", 2)
$d.Content.Find.Execute("879×6=", $true, $false, $false, $false, $false, $true, 1, $false, "314×7=", 2)
$d.Content.Find.Execute("423×4=", $true, $false, $false, $false, $false, $true, 1, $false, "709×6=", 2)
$d.Content.Find.Execute("529×2=", $true, $false, $false, $false, $false, $true, 1, $false, "796×8=", 2)
$d.Content.Find.Execute("697×5=", $true, $false, $false, $false, $false, $true, 1, $false, "475×8=", 2)
$d.Content.Find.Execute("299×4=", $true, $false, $false, $false, $false, $true, 1, $false, "227×5=", 2)
$d.Content.Find.Execute("664×7=", $true, $false, $false, $false, $false, $true, 1, $false, "512×7=", 2)
$d.Content.Find.Execute("336×5=", $true, $false, $false, $false, $false, $true, 1, $false, "315×7=", 2)
$d.Content.Find.Execute("924×4=", $true, $false, $false, $false, $false, $true, 1, $false, "176×3=", 2)
$d.Content.Find.Execute("836×8=", $true, $false, $false, $false, $false, $true, 1, $false, "497×5=", 2)
$d.Content.Find.Execute("101×6=", $true, $false, $false, $false, $false, $true, 1, $false, "735×5=", 2)
$d.Content.Find.Execute("973×6=", $true, $false, $false, $false, $false, $true, 1, $false, "451×9=", 2)
$d.Content.Find.Execute("333×7=", $true, $false, $false, $false, $false, $true, 1, $false, "961×3=", 2)
$d.Content.Find.Execute("376×7=", $true, $false, $false, $false, $false, $true, 1, $false, "233×6=", 2)
$d.Content.Find.Execute("589×2=", $true, $false, $false, $false, $false, $true, 1, $false, "719×8=", 2)
$d.Content.Find.Execute("486×6=", $true, $false, $false, $false, $false, $true, 1, $false, "663×4=", 2)
$d.Content.Find.Execute("888×7=", $true, $false, $false, $false, $false, $true, 1, $false, "575×4=", 2)
$d.Content.Find.Execute("925×4=", $true, $false, $false, $false, $false, $true, 1, $false, "235×4=", 2)
$d.Content.Find.Execute("652×6=", $true, $false, $false, $false, $false, $true, 1, $false, "226×2=", 2)
$d.Content.Find.Execute("967×6=", $true, $false, $false, $false, $false, $true, 1, $false, "708×5=", 2)
$d.Content.Find.Execute("365×5=", $true, $false, $false, $false, $false, $true, 1, $false, "181×4=", 2)
$d.Content.Find.Execute("991×2=", $true, $false, $false, $false, $false, $true, 1, $false, "409×7=", 2)
